# feat: add 2022-Q1 data
#
# Before: "2021-Q3", "2021-Q4", "总计" (totals sheet, 3 rows of data).
# After:  "2021-Q3", "2021-Q4", "2022-Q1" (new fund-level sheet, reusing
#         the old "总计" sheet's id/position), "总计" (new totals sheet,
#         recomputed with the 2022-Q1 quarter folded in as the new first
#         data row).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: duplicate the existing "总计" sheet (3rd tab) *before* touching
# its contents, so the duplicate inherits matching sheetPr/sheetFormatPr/
# pageMargins boilerplate. The duplicate becomes the new "总计" sheet;
# the original (now repurposed) becomes "2022-Q1".
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item(3)
$q1.Name = "2022-Q1"

$q1.Copy($null, $q1)
$total = $wb.Worksheets.Item(4)
$total.Name = "总计"

# ---------------------------------------------------------------------
# Step 2: rebuild "2022-Q1" with fund-level data (same column layout as
# the "2021-Q3"/"2021-Q4" sheets: 基金代码/基金名称/基金规模/股票总仓位/
# 仓位占比/持有市值(亿元)/仓位排名).
# ---------------------------------------------------------------------

# Header style donor: B1 already carries the bold/bordered header style
# used throughout this workbook for header rows, and A2 for the leading
# index column.
$headerStyleCell = $q1.Cells.Item(1, 2)
$indexStyleCell = $q1.Cells.Item(2, 1)

# Re-stamp header row B1:D1 (values change) and extend it with new
# E1:H1 header cells (same header style as B1).
$headerStyleCell.Copy($q1.Cells.Item(1, 3))
$headerStyleCell.Copy($q1.Cells.Item(1, 4))
$headerStyleCell.Copy($q1.Cells.Item(1, 5))
$headerStyleCell.Copy($q1.Cells.Item(1, 6))
$headerStyleCell.Copy($q1.Cells.Item(1, 7))
$headerStyleCell.Copy($q1.Cells.Item(1, 8))

$q1.Cells.Item(1, 2).Value = "基金代码"
$q1.Cells.Item(1, 3).Value = "基金名称"
$q1.Cells.Item(1, 4).Value = "基金规模"
$q1.Cells.Item(1, 5).Value = "股票总仓位"
$q1.Cells.Item(1, 6).Value = "仓位占比"
$q1.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q1.Cells.Item(1, 8).Value = "仓位排名"

# Leading index column (A) style for rows 2-4.
$indexStyleCell.Copy($q1.Cells.Item(2, 1))
$indexStyleCell.Copy($q1.Cells.Item(3, 1))
$indexStyleCell.Copy($q1.Cells.Item(4, 1))

$q1.Cells.Item(2, 1).Value = 0
$q1.Cells.Item(3, 1).Value = 1
$q1.Cells.Item(4, 1).Value = 2

# Row 2 - 010570 新沃创新领航混合A
$q1.Cells.Item(2, 2).Value = "'010570"
$q1.Cells.Item(2, 2).Style = "Normal"
$q1.Cells.Item(2, 3).Value = "新沃创新领航混合A"
$q1.Cells.Item(2, 4).Value = "'0.84"
$q1.Cells.Item(2, 4).Style = "Normal"
$q1.Cells.Item(2, 5).Value = "'93.74"
$q1.Cells.Item(2, 5).Style = "Normal"
$q1.Cells.Item(2, 6).Value = "'4.75"
$q1.Cells.Item(2, 6).Style = "Normal"
$q1.Cells.Item(2, 7).Value = "'0.0399"
$q1.Cells.Item(2, 7).Style = "Normal"
$q1.Cells.Item(2, 8).Value = 6

# Row 3 - 010571 新沃创新领航混合C
$q1.Cells.Item(3, 2).Value = "'010571"
$q1.Cells.Item(3, 2).Style = "Normal"
$q1.Cells.Item(3, 3).Value = "新沃创新领航混合C"
$q1.Cells.Item(3, 4).Value = "'0.56"
$q1.Cells.Item(3, 4).Style = "Normal"
$q1.Cells.Item(3, 5).Value = "'93.74"
$q1.Cells.Item(3, 5).Style = "Normal"
$q1.Cells.Item(3, 6).Value = "'4.75"
$q1.Cells.Item(3, 6).Style = "Normal"
$q1.Cells.Item(3, 7).Value = "'0.0266"
$q1.Cells.Item(3, 7).Style = "Normal"
$q1.Cells.Item(3, 8).Value = 6

# Row 4 - 002564 新沃通盈灵活配置混合
$q1.Cells.Item(4, 2).Value = "'002564"
$q1.Cells.Item(4, 2).Style = "Normal"
$q1.Cells.Item(4, 3).Value = "新沃通盈灵活配置混合"
$q1.Cells.Item(4, 4).Value = "'0.12"
$q1.Cells.Item(4, 4).Style = "Normal"
$q1.Cells.Item(4, 5).Value = "'93.74"
$q1.Cells.Item(4, 5).Style = "Normal"
$q1.Cells.Item(4, 6).Value = "'5.06"
$q1.Cells.Item(4, 6).Style = "Normal"
$q1.Cells.Item(4, 7).Value = "'0.0061"
$q1.Cells.Item(4, 7).Style = "Normal"
$q1.Cells.Item(4, 8).Value = 8

# ---------------------------------------------------------------------
# Step 3: rebuild "总计" with the recomputed totals table (2022-Q1 folded
# in as the new first data row, the older quarters shifting down).
# ---------------------------------------------------------------------

# Row 4 is brand new (the old sheet only had rows 1-3) - stamp it with
# the same leading-index-column style as rows 2/3 before writing values.
$total.Cells.Item(2, 1).Copy($total.Cells.Item(4, 1))

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 3
$total.Cells.Item(2, 4).Value = 0.07000000000000001

$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(3, 2).Value = "2021-Q4"
$total.Cells.Item(3, 3).Value = 6
$total.Cells.Item(3, 4).Value = 0.6

$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(4, 2).Value = "2021-Q3"
$total.Cells.Item(4, 3).Value = 3
$total.Cells.Item(4, 4).Value = 1.16

# Restore the original active/selected tab ("2021-Q3", the 1st sheet) -
# adding/copying sheets shifts Excel's selection to the newest tab.
$wb.Worksheets.Item(1).Activate()
